# Add UUIDs to test data (pulled from #93)
#
# - Column A (rows 3-12) held a simple 1..10 row index; replace each with a
#   generated UUID string (pulled from PR #93).
# - Widen column A so the UUIDs are fully visible.
# - Leave the active selection on the newly-populated UUID column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Patient ID column (A) -> UUIDs instead of 1..10
$uuids = @(
    "fd39991e-8791-4a17-91de-5fee78236d6d",
    "f52d82f3-1fb4-4734-b0d3-56665153243c",
    "d0ffda4a-6e60-4bff-8929-1c5b998c8d28",
    "d0e2eef7-2fc2-41c1-b1a9-226a9686772d",
    "e7defaff-ac0d-409d-a497-cebe94c4ce01",
    "2c8e8f01-7e2b-4ff6-9be6-60fa4914cd24",
    "2204c553-1098-42e6-87cd-8baeedfed672",
    "651bb92d-74af-44d2-a3ca-ef3cbc70ee0c",
    "34c9a874-1227-44d2-9274-3ff757e0fcff",
    "4fd057b7-7078-4900-825c-28d145afa25e"
)

for ($i = 0; $i -lt $uuids.Length; $i++) {
    $row = 3 + $i
    $ws.Range("A$row").Value = $uuids[$i]
}

# Widen column A to comfortably fit the UUID strings.
$ws.Range("A1").EntireColumn.ColumnWidth = 46

# Match the resulting selection left behind in the workbook.
$ws.Range("A3:A12").Select()
